# Update TPM-derived values in the LR-pairs sheet (Gdf2-Acvrl1) per new TPM script output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 31.40242733333333
$ws.Range("N2").Value = 94.20728199999999
$ws.Range("O2").Value = 0.5334014788811394
$ws.Range("P2").Value = 0.5334014788811395
$ws.Range("Q2").Value = 10.55383245294444
$ws.Range("R2").Value = 94.98449207649998
$ws.Range("S2").Value = 0.5334014788811394
$ws.Range("T2").Value = 0.5334014788811395

# Row 3 (only specificity columns change)
$ws.Range("O3").Value = 0.327656036225058
$ws.Range("P3").Value = 0.327656036225058
$ws.Range("S3").Value = 0.327656036225058
$ws.Range("T3").Value = 0.327656036225058

# Row 4
$ws.Range("M4").Value = 1.868202333333333
$ws.Range("N4").Value = 5.604607
$ws.Range("O4").Value = 0.03173327580290011
$ws.Range("P4").Value = 0.03173327580290011
$ws.Range("Q4").Value = 0.6278716675277777
$ws.Range("R4").Value = 5.650845007749999
$ws.Range("S4").Value = 0.03173327580290011
$ws.Range("T4").Value = 0.03173327580290011

# Row 5
$ws.Range("M5").Value = 6.311623666666667
$ws.Range("N5").Value = 18.934871
$ws.Range("O5").Value = 0.1072092090909023
$ws.Range("P5").Value = 0.1072092090909024
$ws.Range("Q5").Value = 2.121231520638889
$ws.Range("R5").Value = 19.09108368575
$ws.Range("S5").Value = 0.1072092090909023
$ws.Range("T5").Value = 0.1072092090909024
